# Add a new "Notre solution" column (G) to the first results table
# (rows 2-29) in the "Feuil1" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, bold like the other header cells in row 2.
$ws.Range("G2").Value = "Notre solution"
$ws.Range("G2").Font.Bold = $true

# Data values for the new column.
$ws.Range("G3").Value = 438
$ws.Range("G4").Value = 642
$ws.Range("G5").Value = 732
$ws.Range("G6").Value = 294
$ws.Range("G7").Value = 750
$ws.Range("G8").Value = 498
$ws.Range("G9").Value = 390
$ws.Range("G10").Value = 522
$ws.Range("G11").Value = 336
$ws.Range("G12").Value = 384
$ws.Range("G13").Value = 258
$ws.Range("G14").Value = 800
$ws.Range("G15").Value = 78
$ws.Range("G16").Value = 136
$ws.Range("G17").Value = 58
$ws.Range("G18").Value = 86
$ws.Range("G19").Value = 131
$ws.Range("G20").Value = 108
$ws.Range("G21").Value = 100
$ws.Range("G22").Value = 97
$ws.Range("G23").Value = 150
$ws.Range("G24").Value = 429
$ws.Range("G25").Value = 453
$ws.Range("G26").Value = 564
$ws.Range("G27").Value = 366
$ws.Range("G28").Value = 328
$ws.Range("G29").Value = 647

# Restore the selection to where the author left it.
$ws.Range("G30").Select() | Out-Null
